# "Minor changes and Interventions runned in large scale mode"
# Switch the Functional Unit (main!C17, "Number of machines") from the
# small pilot-scale figure to the large-scale rollout figure, fix a typo
# in the B33 label, and leave the selection on the label that was edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$ws.Activate()

# Run the intervention in large-scale mode (number of machines).
$ws.Range("C17").Value = 606

# Fix typo: "oter" -> "other".
$ws.Range("B33").Value = "Increase in the use of other services by High Rainfall"

$ws.Range("B33").Select()
